$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Document Style")

# Update the color name for Kenzi's row from "Magenta" to "Red"
$ws.Range("A6").Value = "Red"

# Update the row's fill color from magenta (FFFF66FF) to red (FFFF0000)
$ws.Range("A6:E6").Interior.Color = 255

$ws.Range("C4").Select()
